# CRIANCA_HOSP.xlsx - "All kinds of changes from field tests"
#
# Renames the "region_csv" query/choice-list to "reg_csv", renames the
# "dontknow" choice list used for facility selection to "dontknowfac",
# adds a new "8888" ("Other place") choice alongside the existing "9999"
# ("Don't know") choice, switches the ONDEINTC "other" branch condition
# from '9999' to '8888', and turns the assign-calculation on row 20 of
# the survey sheet into a lookup of the new ondeintcns answer instead of
# a hard-coded 9999 literal.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# survey sheet
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# row 16: values_list region_csv -> reg_csv
$survey.Range("E16").Value = "reg_csv"

# row 18: values_list dontknow -> dontknowfac
$survey.Range("E18").Value = "dontknowfac"

# row 20: calculation literal 9999 -> data('ondeintcns')
$survey.Range("I20").Value = "data('ondeintcns')"

# row 23: condition '9999' -> '8888'
$survey.Range("C23").Value = "data('ONDEINTC') == '8888'"

$survey.Range("E18").Select()

# ---------------------------------------------------------------------
# queries sheet
# ---------------------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")

# row 2: query_name region_csv -> reg_csv
$queries.Range("A2").Value = "reg_csv"

$queries.Range("A2").Select()

# ---------------------------------------------------------------------
# model sheet (selection only moves; no content changes here)
# ---------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")
$model.Activate()
$model.Range("A24").Select()

# ---------------------------------------------------------------------
# choices sheet
# ---------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

# row 9 becomes the new "8888 / Other place" choice ...
$choices.Range("A9").Value = "dontknowfac"
$choices.Range("B9").Formula = "=""8888"""
$choices.Range("C9").Value = "Other place"
$choices.Range("D9").Value = "Outro lugar"

# ... and the old "9999 / Don't know" choice moves down to row 10
$choices.Range("A10").Value = "dontknowfac"
$choices.Range("B10").Formula = "=""9999"""
$choices.Range("C10").Value = "Don't know"
$choices.Range("D10").Value = "Não sabe"

# choices becomes the workbook's active tab, matching the authored edit
$choices.Activate()
$choices.Range("A9:D10").Select()
